$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("datadetails")

$ws.Range("B2").Value = "https://ttri.epixel.link/en/register/"
$ws.Range("B3").Value = "admin"
$ws.Range("B6").Value = "ey093134"
$ws.Range("B7").Value = "abhish097@mailinator.com"
$ws.Range("B11").Value = "akloirt677"

# B12 ("Phone Number") is a digit-only string that must stay text (matches
# the source t="str" cell), not be coerced into a number, so force the
# text number format before assigning the value, then drop back to the
# default style.
$phone = $ws.Range("B12")
$phone.NumberFormat = "@"
$phone.Value = "4347870145"
$phone.Style = "Normal"
